$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in H1, matching the style of the other headers (B1:G1)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in Label column values for the first block (rows 2-11)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# Fill in Label column values for the second block (rows 12-21)
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1

# Update refit values for Control 13 / Control 51 / MDD 50 at 100 iterations
$ws.Range("D4").Value = 0.4377385299476793
$ws.Range("E4").Value = 0.4377385299476793

$ws.Range("D6").Value = 0.6229160857056946
$ws.Range("E6").Value = 0.6229160857056946

$ws.Range("D9").Value = 0.4091159723620472
$ws.Range("E9").Value = 0.5908840276379528
